# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two middle workers (DIANA CAROLINA ARIAS ROMERO, LUZ DARY MARQUEZ GUERRERO)
# which occupied rows 17-19; this leaves the remaining MARIA JOSE rows (old 20-23) shifted
# up into rows 17-20, correctly carrying the bottom-border row style up with them.
$ws.Range("A17:A19").EntireRow.Delete()

# Update the account-summary figures
$ws.Range("E11").Value = 345364
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 5

# Row 16 (GUILLERMO ALBERTO OLARTE GARCIA, periodo 2002) - valor mora updated
$ws.Range("G16").Value = 5000000

# Rows 17-20 now all belong to MARIA JOSE PADILLA CASTILLO; fix the periodo order
# to ascending (2111, 2112, 2201, 2202) as in the refreshed source data.
$ws.Range("E17").Value = "2111"
$ws.Range("E18").Value = "2112"
$ws.Range("E19").Value = "2201"
$ws.Range("E20").Value = "2202"

$wb.Save()
